$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.434.64'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.305.86'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.36'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.44'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.82%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.303.73'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.887.90'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '66.408.86'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.66'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.320.68'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '423.21'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.12'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.32'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.08%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.42'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.57%  '
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.67'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.465.18'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.515'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.19%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.205'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.79%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000115'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.29%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.11'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.39%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.92'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.36'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.16'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.36%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.59'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.50%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.18'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.51%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.04'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.43'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.80'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.856.17'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.27'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.23%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.33'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.753'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.62%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.73'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0658'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.89'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.55%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.07'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.52%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '311.19'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -4.37%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0272'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.16%  '
